# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The workbook originally listed 3 worker-period records (rows 16-18) for
# the account holder. This edit collapses it down to a single record
# (the 2211-period / 40000 mora row that used to be row 17), updates the
# summary totals accordingly, and removes the now-obsolete worker rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Summary block: total mora owed, worker count, period count.
$ws.Range("E11").Value = 40000   # VALOR MORA total: 248940 -> 40000
$ws.Range("C13").Value = 1       # Cant. Trabajadores: 2 -> 1
$ws.Range("F13").Value = 1       # Cant. Periodos: 2 -> 1

# Row 16 becomes the surviving detail record: period 2211, mora 40000
# (previously this row held the 2507 / 152000 record).
$ws.Range("E16").Value = "2211"
$ws.Range("F16").Value = 40000

# Remove the two obsolete worker/period rows (old row 17: 2211 duplicate,
# old row 18: OSCAR DAVID MARTINEZ GARCIA). Excel shifts everything below
# up automatically, so the signature block (old rows 23/24) becomes 21/22.
$ws.Range("17:18").EntireRow.Delete()
